$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.991.24"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "1.900.37"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.7869"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "'244.51"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'25.90"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").Value = "'0.07323"
$ws.Range("E10").Value = "  +4.92%  "
$ws.Range("D11").Value = "'0.08115"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "'0.7788"
$ws.Range("E12").Value = "  +3.11%  "
$ws.Range("D13").Value = "'5.518"
$ws.Range("E13").Value = "  +4.62%  "
$ws.Range("D14").Value = "1.879.81"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").Value = "'94.51"
$ws.Range("E15").Value = "  +2.67%  "
$ws.Range("D16").Value = "'6.270"
$ws.Range("E16").Value = "  +6.19%  "
$ws.Range("D17").Value = "29.945.72"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").Value = "'246.76"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").Value = "'0.000007843"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").Value = "'8.175"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D23").Value = "2.120.14"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").Value = "'1.0000"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'0.1603"
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("D26").Value = "'9.495"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("D27").Value = "'163.32"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "'2.043"
$ws.Range("D30").Value = "'1.437"
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").Value = "'1.548"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "'4.495"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D33").Value = "'0.05614"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").Value = "'1.253"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "'0.7565"
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("D37").Value = "'0.9997"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'2.678"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").Value = "'2.790"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").Value = "1.148.71"
$ws.Range("E41").Value = "  +13.05%  "
$ws.Range("D42").Value = "'0.4479"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("D43").Value = "'74.31"
$ws.Range("E43").Value = "  +3.04%  "
$ws.Range("D44").Value = "'5.989"
$ws.Range("E44").Value = "  +3.30%  "
$ws.Range("D45").Value = "'0.8564"
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("D46").Value = "'1.905"
$ws.Range("E46").Value = "  +2.46%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "'3.162"
$ws.Range("E48").Value = "  +9.14%  "
$ws.Range("D49").Value = "'102.27"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'7.561"
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.799"
$ws.Range("E51").Value = "  -0.74%  "
